$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("PCA-300-Corpus")
$dst = $wb.Worksheets.Item("PCA_100_corpus")

# Seed the new sheet with the same layout/formatting (borders, fonts,
# wrap-text, shared-string reuse) as the already-populated "PCA-300-Corpus"
# results table, then correct the header labels and numbers that differ.
$src.Range("A1:M4").Copy($dst.Range("A1"))

# Header row (row 1) - fix the few headers that differ from the source sheet
$dst.Range("A1").Value = "ngram"
$dst.Range("C1").Value = " logreg_test"
$dst.Range("F1").Value = "svm_train"
$dst.Range("J1").Value = "adaboost_train"
$dst.Range("L1").Value = "rf_train"

# Row 2 data (1gram)
$dst.Range("B2").Value = 0.64399300000000004
$dst.Range("C2").Value = 0.63981600000000005
$dst.Range("D2").Value = 0.54246799999999995
$dst.Range("E2").Value = 0.53738900000000001
$dst.Range("F2").Value = 0.63775599999999999
$dst.Range("G2").Value = 0.63446000000000002
$dst.Range("H2").Value = 0.55946399999999996
$dst.Range("I2").Value = 0.558083
$dst.Range("J2").Value = 0.68155699999999997
$dst.Range("K2").Value = 0.65675399999999995
$dst.Range("L2").Value = 0.59143900000000005
$dst.Range("M2").Value = 0.58886300000000003

# Row 3 data (2gram)
$dst.Range("B3").Value = 0.52627100000000004
$dst.Range("C3").Value = 0.52681599999999995
$dst.Range("D3").Value = 0.50060300000000002
$dst.Range("E3").Value = 0.50351299999999999
$dst.Range("F3").Value = 0.524509
$dst.Range("G3").Value = 0.52392899999999998
$dst.Range("H3").Value = 0.53748200000000002
$dst.Range("I3").Value = 0.52980700000000003
$dst.Range("J3").Value = 0.57601999999999998
$dst.Range("K3").Value = 0.55585700000000005
$dst.Range("L3").Value = 0.53547699999999998
$dst.Range("M3").Value = 0.53119799999999995

# Row 4 data (3gram)
$dst.Range("B4").Value = 0.469057
$dst.Range("C4").Value = 0.468802
$dst.Range("D4").Value = 0.45826299999999998
$dst.Range("E4").Value = 0.45972499999999999
$dst.Range("F4").Value = 0.46887099999999998
$dst.Range("G4").Value = 0.468663
$dst.Range("H4").Value = 0.47992000000000001
$dst.Range("I4").Value = 0.47808800000000001
$dst.Range("J4").Value = 0.49047000000000002
$dst.Range("K4").Value = 0.48100999999999999
$dst.Range("L4").Value = 0.47603600000000001
$dst.Range("M4").Value = 0.47579300000000002

# Row heights to match the other results tables
$dst.Rows.Item(1).RowHeight = 46
$dst.Rows.Item(2).RowHeight = 17
$dst.Rows.Item(3).RowHeight = 17
$dst.Rows.Item(4).RowHeight = 17

# Make "PCA_100_corpus" the active sheet/tab and select C10 on it, which also
# clears the previous tab-selected flag on whichever sheet had it before.
$dst.Activate()
$dst.Range("C10").Select()

Write-Host "Done"
